$d = $word.ActiveDocument
$wdParagraph = 4

# ---------------------------------------------------------------------------
# 1) "water_shortage_level" -> "actual_water_shortage_level" (bold heading,
#    stays inside the existing spell-check markers, keeps bold/bCs).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("water_shortage_level", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "actual_water_shortage_level", 2) | Out-Null
$afterHeading = $rng.End

# ---------------------------------------------------------------------------
# 2) Fill in the placeholder "This table summarizes…" paragraph with the
#    full blurb. The paragraph already owns <w:commentRangeStart w:id="0"/>
#    (the host preserves non-run children automatically), so the new XML
#    only needs to carry the replacement runs - including a
#    <w:lastRenderedPageBreak/> landing in front of "Reporting.".
# ---------------------------------------------------------------------------
$rng = $d.Range($afterHeading, $d.Content.End)
$rng.Find.ClearFormatting()
$rng.Find.Execute("This table summarizes" + [char]0x2026) | Out-Null
$blurbPara = $rng.Duplicate
$blurbPara.Expand($wdParagraph)
$afterBlurbSearch = $rng.Start

$apos = [char]0x2019
$blurbXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r><w:t xml:space="preserve">This table </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">reports the monthly state standard shortage level </w:t></w:r>' + `
  '<w:r><w:t>by urban retail water suppliers, which are generally defined as agencies serving over 3,000 service connections or deliveries 3,000 acre-feet of water annually for municipal purposes.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> These data</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> are collected by the State Water Resources Control Board</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> t</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">hrough </w:t></w:r>' + `
  ('<w:r><w:t xml:space="preserve">it' + $apos + 's monthly Conservation </w:t></w:r>') + `
  '<w:r><w:lastRenderedPageBreak/><w:t>Reporting.</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> The data reported</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> for this project have been filtered to show the most recent data (M</w:t></w:r>' + `
  '<w:r><w:t>arch 2022 and on).</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '</w:p>'

$blurbPara.InsertXML($blurbXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Add an italic " " and a plain "?" right after the comment reference in
#    the "Data use limitations:" paragraph that follows the blurb. The
#    paragraph already owns the <w:commentRangeEnd w:id="0"/> and the
#    <w:commentReference w:id="0"/> run, so only the two appended runs need
#    to be supplied.
# ---------------------------------------------------------------------------
$rng = $d.Range($afterHeading, $d.Content.End)
$rng.Find.ClearFormatting()
$rng.Find.Execute("Data use limitations:") | Out-Null
$limPara = $rng.Duplicate
$limPara.Expand($wdParagraph)

$limXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Data use limitations:</w:t></w:r>' + `
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t>?</w:t></w:r>' + `
  '</w:p>'

$limPara.InsertXML($limXml) | Out-Null

# ---------------------------------------------------------------------------
# 4) Drop the now-stale <w:lastRenderedPageBreak/> in front of
#    "Other Relevant " - the (re-)rendered page break now lands earlier,
#    inside the blurb paragraph above. The paragraph keeps its Heading2
#    pPr automatically, so the new XML only needs the pStyle + runs.
# ---------------------------------------------------------------------------
$rng = $d.Range($afterHeading, $d.Content.End)
$rng.Find.ClearFormatting()
$rng.Find.Execute("Other Relevant ") | Out-Null
$headingPara = $rng.Duplicate
$headingPara.Expand($wdParagraph)

$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r><w:t xml:space="preserve">Other Relevant </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">California State </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">Open Data </w:t></w:r>' + `
  '</w:p>'

$headingPara.InsertXML($headingXml) | Out-Null

Write-Host "Edits applied."
